$d = $word.ActiveDocument

# --- Update existing paragraphs 1-6 (full text replacement) ---
$d.Paragraphs.Item(1).Range.Text = "⚡️🚀המאמר היומי של מייק 03.07.24:⚡️🚀"
$d.Paragraphs.Item(2).Range.Text = "The Remarkable Robustness of LLMs: Stages of Inference?"
$d.Paragraphs.Item(3).Range.Text = "מאמר מעניין החוקר איזה שכבות ניתן לזרוק ממודל השפה ועדיין לשמור על ביצועים נאותים. אתם אולי מכירים lottery ticket hypothesis הטוען כי ברשתות עתירות פרמטרים (overparameterized) בד:כ ניתן למצוא קטנה הרבה יותר עם ביצועים מאוד קרובים אך הבעיה שאנו לא יודעים לאתר אותה."
$d.Paragraphs.Item(4).Range.Text = "המאמר כאמור בחן איזה שכבות הן סוג של מיותרות במודלי שפה והגיע לתופעות מעניינות לגבי תהליך האינפרנס שלהם. הם זיהו 4 שלבים עיקריים"
$d.Paragraphs.Item(5).Range.Text = "דה-טוקניזציה או רכישה התחלתית של קשרים קונטקסטואליים: טרנספורמציה ראשונית של ייצוג ה-raw (מהמילון) של הטוקנים לייצוג תלוי הקשר (חישובי attention כבדים לכל אורך הקונטקסט)."
$d.Paragraphs.Item(6).Range.Text = "הנדסת פיצ'רים התחלתיים מהייצוגים תלוי הקשר מהשלב הקודם ו״הכנת קרקע״ לחיזוי של הטוקנים הבאים. עדיין לא ניתן לחזות את הטוקנים האלו מהפיצ'רים בשלב הזה אבל המודל מתחיל ״להבין הקשרים מרחבים ועתיים בטקסט (היה מחקר מעניין הזה)"

# --- Paragraph 7: was the arXiv link, now becomes new body paragraph ---
$d.Paragraphs.Item(7).Range.Text = "בניית קבוצות נוירונים (אנסמבל) לחיזוי הטוקן הבא. בשלב הזה הרשת מתחילה להתכנס ולבנות קבוצות ״prediction neurons`" שישולבו יחד למטרת חיזוי הטוקן הבא."

# --- Append new trailing paragraphs (8-11), preserving Normal style ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "חידוד של prediction neurons: הרשת ״בוחרת״ את הנוירונים החשובים ביותר לחיזוי הטוקן הבא על ידי הדעכה של חלק מה-prediction neurons מהשלב הקודם."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "והכי חשוב שהשכבות מעורבות בשלב 1 ובשלב 4 הם הכי חשובות לביצוע המודל כאשר חלק מהשכבות של שלב 2 ו-3 ניתן להסיר ללא פגיעה משמעותית בביצועים. "
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "הרבה טענות מעניינות במאמר הזה (חלקם הגדול זה סיכום של העבודות הקודמות בנושא הזה)."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "https://arxiv.org/abs/2406.19384"

Write-Output "edit applied"